$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update Runmode column (C) from YES to NO for rows 2,4,5,6,7
$ws.Range("C2").Value = "NO"
$ws.Range("C4").Value = "NO"
$ws.Range("C5").Value = "NO"
$ws.Range("C6").Value = "NO"
$ws.Range("C7").Value = "NO"

# Add new rows 8 and 9
$ws.Range("A8").Value = "Profile"
$ws.Range("B8").Value = "Description"
$ws.Range("C8").Value = "YES"

$ws.Range("B9").Value = "User Dashboard Scenarios"
$ws.Range("A9").Value = "UserDashboard"
$ws.Range("C9").Value = "YES"

# Update the selection to match the post-edit state
$ws.Range("B9").Select()
